$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 43, shifting existing rows 43:88 down to 44:89
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record
$ws.Cells.Item(43, 1).Value = 7
$ws.Cells.Item(43, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(43, 3).Value = "Ñuble"
$ws.Cells.Item(43, 4).Value = 44763
$ws.Cells.Item(43, 5).Value = 16
$ws.Cells.Item(43, 6).Value = 100112021
$ws.Cells.Item(43, 7).Value = "Ají"
$ws.Cells.Item(43, 8).Value = "Inferno"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 40
$ws.Cells.Item(43, 11).Value = 18000
$ws.Cells.Item(43, 12).Value = 18000
$ws.Cells.Item(43, 13).Value = 18000
$ws.Cells.Item(43, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 1200
$ws.Cells.Item(43, 17).Value = 15
$ws.Cells.Item(43, 18).Value = "Hortaliza"
